$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (prices) that must remain
# literal text (e.g. "30.667.42", "1.000", "0.000007926"). Force the
# range to Text format before assigning so Excel does not coerce the
# string into a Number/Double and alter its representation.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.667.42'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.891.27'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.90'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4965'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2965'
$ws.Range('E8').Value = '  +2.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06816'
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.891.13'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.07'
$ws.Range('E11').Value = '  +2.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07314'
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '90.79'
$ws.Range('E13').Value = '  +5.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.077'
$ws.Range('E14').Value = '  +4.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6732'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.657.11'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007926'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.21'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.134.77'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.859'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '176.42'
$ws.Range('E23').Value = '  +31.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.062'
$ws.Range('E24').Value = '  +8.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.271'
$ws.Range('E25').Value = '  +2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.58'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.75'
$ws.Range('E27').Value = '  +12.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.927'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.392'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.337'
$ws.Range('E30').Value = '  +4.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08929'
$ws.Range('E31').Value = '  +2.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.030'
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05228'
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7415'
$ws.Range('E34').Value = '  +5.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.138'
$ws.Range('E35').Value = '  +3.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.676'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01879'
$ws.Range('E37').Value = '  +10.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.701'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.172'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9357'
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4372'
$ws.Range('E41').Value = '  +4.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.36'
$ws.Range('E42').Value = '  +3.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.814'
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.677'
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('E46').Value = '  +8.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05841'
$ws.Range('E47').Value = '  +3.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.43'
$ws.Range('E48').Value = '  +3.28%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3880'
$ws.Range('E49').Value = '  +5.34%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.497'
$ws.Range('E50').Value = '  +4.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.381'
$ws.Range('E51').Value = '  +3.75%  '
